# "Add files via upload" -- populate the "17" submission-group worksheet
# (the 7th sheet) with the newly uploaded names, and make it the active
# sheet/selection, matching the other populated group sheets ("12", "13").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("17")

$names = @("Yuval Koskas", "Lior Tsalovich", "Noam Raanan", "Maxim Gutnik")

for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $names[$i]
}

$ws.Activate()
$ws.Range("D9").Select()
